$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the test row (row 5: "Luis" / "XTE - XT & Tools")
$ws.Range("A5:T5").EntireRow.Delete()

# Remove the test column (column T: "Teste")
$ws.Range("T1:T4").EntireColumn.Delete()

# Uncheck D4 (was TRUE, now FALSE)
$ws.Cells.Item(4, 4).Value = $false
